$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are kept as Text (matches original inlineStr semantics),
# since Excel would otherwise auto-coerce numeric-looking strings (e.g. "10.10" -> 10.1).
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D10", "D11", "D13", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D24", "D28", "D32", "D33", "D34", "D35", "D36", "D38", "D42", "D43", "D47", "D48", "D49", "D50", "D51")
foreach ($pc in $priceCells) { $ws.Range($pc).NumberFormat = "@" }

$ws.Range("D2").Value = "67.381.75"
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("D3").Value = "3.704.84"
$ws.Range("E3").Value = "  -4.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "596.46"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("D6").Value = "165.53"
$ws.Range("E6").Value = "  -5.69%  "
$ws.Range("D7").Value = "3.702.97"
$ws.Range("E7").Value = "  -4.34%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").Value = "  -4.82%  "
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("D13").Value = "37.61"
$ws.Range("E13").Value = "  -6.24%  "
$ws.Range("E14").Value = "  -5.40%  "
$ws.Range("D15").Value = "4.326.44"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "3.708.68"
$ws.Range("E16").Value = "  -4.18%  "
$ws.Range("D17").Value = "67.456.50"
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("E18").Value = "  +5.23%  "
$ws.Range("D19").Value = "7.18"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").Value = "0.114"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "486.92"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("D22").Value = "9.44"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("D24").Value = "85.37"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").Value = "  -6.84%  "
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("E27").Value = "  -4.00%  "
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  -7.61%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.61"
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "31.41"
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").Value = "3.845.00"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.648.54"
$ws.Range("E35").Value = "  -4.23%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.107"
$ws.Range("E36").Value = "  -5.25%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("E39").Value = "  -6.43%  "
$ws.Range("E40").Value = "  -7.74%  "
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("D42").Value = "432.87"
$ws.Range("E42").Value = "  -9.32%  "
$ws.Range("D43").Value = "48.58"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("E45").Value = "  -6.76%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "40.66"
$ws.Range("E48").Value = "  -6.35%  "
$ws.Range("D49").Value = "142.34"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "2.753.19"
$ws.Range("E50").Value = "  -6.29%  "
$ws.Range("D51").Value = "0.0348"
$ws.Range("E51").Value = "  -4.07%  "
